$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, copying the formatting used by the other
# header cells (B1:G1) so it matches the bold/centered/bordered style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = 0
$ws.Range("H1").Value = "Save"

# Fill H2:H21 with 0 values (new "Save" data column)
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
